$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column P (16) is widened slightly for readability ("app friendliness").
#    The engine quantizes column widths to whole-pixel steps, so we drive
#    ColumnWidth to the value that lands on the nearest achievable step.
# ---------------------------------------------------------------------------
$ws.Columns.Item(16).ColumnWidth = 11.6

# ---------------------------------------------------------------------------
# 2. Fix orientation/text issues: the label columns -- Subject (A),
#    Scan_Date (B), Process_Date (E) and Healthy_Cohort (F) -- plus the
#    header row, are explicitly stored as Text ("@") so date-like strings
#    such as "2021-11-19" are never silently reinterpreted as date serials.
# ---------------------------------------------------------------------------
$ws.Range("A1:BN1").NumberFormat = "@"
$ws.Range("A2:B4").NumberFormat = "@"
$ws.Range("E2:F4").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 3. Append two new result rows (5 and 6): subject Xe-020, scanned
#    2021-08-27, reprocessed on 2021-11-19. Text formatting on the label
#    columns is applied before the values are written so the Process_Date
#    string is kept as text instead of being parsed as a date.
# ---------------------------------------------------------------------------
$ws.Range("A5:B6").NumberFormat = "@"
$ws.Range("E5:F6").NumberFormat = "@"

$newRowValues = @(
    "Xe-020", "2021-08-27", 0.46000000000000002, 14.6, "2021-11-19", "20210218_FullHealthyCohort",
    -7399.5510609391858, -663.52565132208792, -2.3745468730668358, 10.587957922601987, 1.6488518829155827, 1.0034600165302596,
    1.2789117068403721, -0.16188560655732437, 82.927229389685365, 170.71287741585741, 0.29767026936348273, 0.10475470466717995,
    105.28697933648426, 36.639497232845422, 83.009067352349163, 42.014908816323285, 28.283825254705842, 8.1450858006302198,
    0.42553056237505466, 0.23300773425490726, 0.68276174457409, 0.21774583292331645, 0.6377339981515634, 0.17700493332617426,
    0.20570714326443451, 0.16870417576884322, 0.33176003454116748, 0.18317431032862019, 10.009671179883945, 11.738394584139265,
    25.386847195357831, 32.761121856866538, 16.477272727272727, 3.6266924564796903, 4.1644277270284791, 39.91133799032778,
    48.992477162815689, 6.4615797958087056, 0.38957549704459971, 0.080601826974744759, 0.99408919935518536, 21.010209564750134,
    62.305212251477705, 14.212788823213327, 1.370231058570661, 0.053734551316496508, 0, 0.053734551316496508,
    11.055883933369156, 45.09672219236969, 41.026329930145081, 2.7001612036539493, 0.06716818914562063, 0.053734551316496508,
    11.069317571198281, 25.819451907576575, 54.083825900053739, 8.4497581945190756, 0.53734551316496504, 0.040300913487372379
)

for ($col = 1; $col -le 66; $col++) {
    $value = $newRowValues[$col - 1]
    $ws.Cells.Item(5, $col).Value = $value
    $ws.Cells.Item(6, $col).Value = $value
}
